# The deck's "datetimeFigureOut" date field (shown on every slide layout
# and on the slide master via the Date placeholder) was refreshed from
# 12/5/2023 to 1/26/2024. Update the Date placeholder's text wherever it
# appears: on the slide master and on every one of its custom layouts.

$p = $ppt.ActivePresentation
$newDate = "1/26/2024"
$ppPlaceholderDate = 16

function Update-DatePlaceholders {
    param($container)

    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)

        if (-not $shp.HasTextFrame) { continue }

        $isDatePlaceholder = $false
        if ($shp.Type -eq 14) {
            try {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
        }

        if ($isDatePlaceholder) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$master = $p.SlideMaster

# Slide master's own Date placeholder.
Update-DatePlaceholders $master

# Every custom (slide) layout hanging off the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout
}
